$d = $word.ActiveDocument

# Locate the paragraph that finishes the 16.04.2021 entry ("...own app.").
# All of the new 17.04.2021 diary content is appended right after it, before
# the two trailing blank paragraphs that close the document body.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Now I think is the time to start planning for my own app.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchor = $findRange.Paragraphs(1)

# Word's AutoCorrect turns a typed straight apostrophe into a curly right
# single quotation mark (U+2019); reproduce that directly since
# Range.InsertAfter does not run AutoCorrect for us.
$rsquo = [char]0x2019
$cr = "`r"

$entry1 = "Starting my own project. I have decided to also use Materializecss to help me build a decent looking site. I am not 100% how I am going to structure this project or where. Guide in the moodle says that I should make an extra folder where I would build my pages. At least in a Git repo, but surely I can" + $rsquo + "t make a new project and push that into a separate folder. Also I am not sure how does the build process  work if my html and js files are not in dist-folder, but in new coursework-folder. I think I need to ponder this a bit! I am going to watch the first and last videos again and see if I can modify some settings."

$entry2 = "I got the project working in s new directory. Decided for scss to use the same folder and mark new files beginning with coursework_.  On my index.html I am linking stylesheet to css/coursework_main.css. Package.json file" + $rsquo + "s sass script I changed to node-sass -w scss/ -o coursework/css/ --recursive. This way everything now seems to point to right directions and the build can commence and the next thing to figure out is the installation of the materialize. Should I npm it or just use CDN."

$entry3 = "With CDN it was very easy to get it working, but then I wanted to change the primary color and it seemed for that the SASS was better choice with variables. Copied the relevant files to correct folders and now I have huge variety of options to choose from."

$block = $cr + $cr + "17.04.2021" + $cr + $entry1 + $cr + $cr + $entry2 + $cr + $cr + $entry3 + $cr

$anchor.Range.InsertAfter($block)
